# New changes - 9/6/20
# Adds two new worksheets ("InvokeviaTxt" and "Sheet2") after "MeetingResume",
# populates them with request data mirroring the existing sheets' layout,
# and moves the active-sheet/selection state to the new last sheet.

$wb = $excel.ActiveWorkbook

$sheet1        = $wb.Worksheets.Item(1)
$meetingResume = $wb.Worksheets.Item("MeetingResume")

# MeetingResume stops being the active tab; its selection reverts to the
# whole-sheet default instead of the old "last edited cell" selection.
[void]$meetingResume.Cells.Select()

# --- New sheet: InvokeviaTxt (sheetId 5), inserted right after MeetingResume ---
$invokeViaTxt = $wb.Worksheets.Add($null, $meetingResume)
$invokeViaTxt.Name = "InvokeviaTxt"

# Seed values + formatting (fills, wrap, vertical alignment) from Sheet1 so
# the new rows line up with the rest of the workbook's styling, then
# overwrite only the cells that actually differ for this request.
$sheet1.Range("A1:E3").Copy($invokeViaTxt.Range("A1:E3"))

$invokeViaTxt.Range("A2").Value = "/cvi/dm/api/v1/invoke/text/json?intent=true&skill=true"
$invokeViaTxt.Range("A3").Value = "/cvi/dm/api/v1/invoke/text/json?intent=true&skill=true"
$invokeViaTxt.Range("E2").Value = "{
  ""text"": ""starte das toronto meeting""
}"
$invokeViaTxt.Range("E3").Value = "{
  ""text"": ""starte das Messe meeting""
}"
$invokeViaTxt.Rows.Item(2).RowHeight = 43.5
$invokeViaTxt.Rows.Item(3).RowHeight = 43.5

[void]$invokeViaTxt.Cells.Select()

# --- New sheet: Sheet2 (sheetId 6), inserted right after InvokeviaTxt ---
$sheet2 = $wb.Worksheets.Add($null, $invokeViaTxt)
$sheet2.Name = "Sheet2"

$sheet1.Range("A1:E3").Copy($sheet2.Range("A1:E3"))

$sheet2.Range("D1").Value = "Xtenant"
$sheet2.Range("A2").Value = "/svhb/meeting/v1/ui/transcript/?meetingId=3"""
$sheet2.Range("D2").Value = "GLOBAL"
$sheet2.Range("A3").Value = "/svhb/meeting/v1/ui/transcript/?meetingId=1"""
$sheet2.Range("D3").Value = "GLOBAL"
# This sheet has no "Body" column - column E stays blank (format only).
$sheet2.Range("E1:E3").ClearContents()

# Final UI state: Sheet2 is the active sheet/tab, with B3 selected.
[void]$sheet2.Range("B3").Select()
